# Auto-generated Excel COM-interop edit script
$wb = $excel.ActiveWorkbook

# --- Step 1: strip the bold/bordered/centered header style (row 1) on every sheet ---
foreach ($ws in $wb.Worksheets) {
    $ws.Range("A1:N1").Style = "Normal"
}

# --- Step 2: update recomputed profit-sheet figures ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 961.6087
$ws.Range("I28").Value = 489.75
$ws.Range("K28").Value = 489.75
$ws.Range("M28").Value = -4.75
$ws.Range("H76").Value = 4999.6665
$ws.Range("I76").Value = 4999
$ws.Range("K76").Value = 4999
$ws.Range("M76").Value = -4684
$ws.Range("H79").Value = 4999.6665
$ws.Range("I79").Value = 4999
$ws.Range("K79").Value = 4999
$ws.Range("M79").Value = -3907
$ws.Range("H86").Value = 7532.8335
$ws.Range("I86").Value = 998.5
$ws.Range("K86").Value = 998.5
$ws.Range("M86").Value = 124.5
$ws.Range("H89").Value = 7532.8335
$ws.Range("I89").Value = 998.5
$ws.Range("K89").Value = 4992.5
$ws.Range("M89").Value = 623.5
$ws.Range("H96").Value = 344
$ws.Range("I96").Value = 344
$ws.Range("K96").Value = 1032
$ws.Range("M96").Value = 341
$ws.Range("H107").Value = 1204.3889
$ws.Range("I107").Value = 1199.9231
$ws.Range("K107").Value = 1199.9231
$ws.Range("M107").Value = 720.0769
$ws.Range("H116").Value = 66892.18
$ws.Range("I116").Value = 149924.86
$ws.Range("J116").Value = 8769.3
$ws.Range("K116").Value = 149924.86
$ws.Range("L116").Value = 8769.3
$ws.Range("M116").Value = -146482.86
$ws.Range("N116").Value = -15653.3
$ws.Range("H132").Value = 62806.53
$ws.Range("J132").Value = 3132.4
$ws.Range("L132").Value = 9397.2
$ws.Range("N132").Value = -14457.2
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3503.138
$ws.Range("I32").Value = 3551.6155
$ws.Range("J32").Value = 3083
$ws.Range("K32").Value = 3551.6155
$ws.Range("L32").Value = 3083
$ws.Range("M32").Value = -3264.6155
$ws.Range("N32").Value = -3657
$ws.Range("H61").Value = 4000.2083
$ws.Range("I61").Value = 3961.1304
$ws.Range("K61").Value = 3961.1304
$ws.Range("M61").Value = -3749.1304
$ws.Range("H63").Value = 2944.4
$ws.Range("I63").Value = 2868.75
$ws.Range("K63").Value = 2868.75
$ws.Range("M63").Value = -2182.75
$ws.Range("H66").Value = 2944.4
$ws.Range("I66").Value = 2868.75
$ws.Range("K66").Value = 14343.75
$ws.Range("M66").Value = -10911.75
$ws.Range("H74").Value = 3133.2122
$ws.Range("I74").Value = 3257.9678
$ws.Range("J74").Value = 1199.5
$ws.Range("K74").Value = 3257.9678
$ws.Range("L74").Value = 1199.5
$ws.Range("M74").Value = -2383.9678
$ws.Range("N74").Value = -2947.5
$ws.Range("H77").Value = 3133.2122
$ws.Range("I77").Value = 3257.9678
$ws.Range("J77").Value = 1199.5
$ws.Range("K77").Value = 16289.839
$ws.Range("L77").Value = 5997.5
$ws.Range("M77").Value = -11921.839
$ws.Range("N77").Value = -14733.5
$ws.Range("H110").Value = 3476.7727
$ws.Range("I110").Value = 2984.5715
$ws.Range("K110").Value = 2984.5715
$ws.Range("M110").Value = -939.5715
$ws.Range("H122").Value = 6474.5
$ws.Range("I122").Value = 6474.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 19423.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -16973.5
$ws.Range("N122").ClearContents()
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
$ws.Range("H136").Value = 4000.2083
$ws.Range("I136").Value = 3961.1304
$ws.Range("K136").Value = 11883.3912
$ws.Range("M136").Value = -9333.3912

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4932.1113
$ws.Range("I105").Value = 1096
$ws.Range("J105").Value = 8001
$ws.Range("K105").Value = 1096
$ws.Range("L105").Value = 8001
$ws.Range("M105").Value = 651
$ws.Range("N105").Value = -11495
$ws.Range("H107").Value = 5737.8857
$ws.Range("I107").Value = 6327.6294
$ws.Range("K107").Value = 6327.6294
$ws.Range("M107").Value = -4407.6294
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H134").Value = 7972.75
$ws.Range("I134").Value = 1899
$ws.Range("J134").Value = 9997.333
$ws.Range("K134").Value = 5697
$ws.Range("L134").Value = 29991.999
$ws.Range("N134").Value = -35061.999
$ws.Range("M134").Value = -3162  # new cell

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 727.5
$ws.Range("I22").Value = 484.25
$ws.Range("J22").Value = 1700.5
$ws.Range("K22").Value = 484.25
$ws.Range("L22").Value = 1700.5
$ws.Range("M22").Value = -134.25
$ws.Range("N22").Value = -2400.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 830.55
$ws.Range("I2").Value = 924
$ws.Range("J2").Value = 550.2
$ws.Range("K2").Value = 5544
$ws.Range("L2").Value = 3301.2
$ws.Range("M2").Value = -5431
$ws.Range("N2").Value = -3527.2
$ws.Range("H40").Value = 67.695656
$ws.Range("I40").Value = 70
$ws.Range("J40").Value = 63.375
$ws.Range("K40").Value = 280
$ws.Range("L40").Value = 253.5
$ws.Range("M40").Value = -211
$ws.Range("N40").Value = -391.5
$ws.Range("H44").Value = 861.8
$ws.Range("I44").Value = 851.75
$ws.Range("J44").Value = 902
$ws.Range("K44").Value = 2555.25
$ws.Range("L44").Value = 2706
$ws.Range("M44").Value = -2157.25
$ws.Range("N44").Value = -3502
$ws.Range("H55").Value = 5472.7144
$ws.Range("J55").Value = 7262
$ws.Range("L55").Value = 21786
$ws.Range("N55").Value = -22140
$ws.Range("H68").Value = 1613.9333
$ws.Range("J68").Value = 1792.5
$ws.Range("L68").Value = 5377.5
$ws.Range("N68").Value = -6999.5
$ws.Range("H71").Value = 1613.9333
$ws.Range("J71").Value = 1792.5
$ws.Range("L71").Value = 16132.5
$ws.Range("N71").Value = -24244.5
$ws.Range("H86").Value = 297.14285
$ws.Range("I86").Value = 297.14285
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 891.4285500000001
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = 294.5714499999999
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 297.14285
$ws.Range("I89").Value = 297.14285
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 2674.28565
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = 3253.71435
$ws.Range("N89").ClearContents()
$ws.Range("H109").Value = 1927.1875
$ws.Range("I109").Value = 1389
$ws.Range("J109").Value = 10000
$ws.Range("K109").Value = 4167
$ws.Range("L109").Value = 30000
$ws.Range("M109").Value = -3127
$ws.Range("N109").Value = -32080
$ws.Range("H113").Value = 517.4167
$ws.Range("I113").Value = 638
$ws.Range("J113").Value = 431.2857
$ws.Range("K113").Value = 1914
$ws.Range("L113").Value = 1293.8571
$ws.Range("M113").Value = 256
$ws.Range("N113").Value = -5633.8571

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8054.6665
$ws.Range("I70").Value = 10374.5
$ws.Range("K70").Value = 10374.5
$ws.Range("M70").Value = -10104.5
$ws.Range("H73").Value = 8054.6665
$ws.Range("I73").Value = 10374.5
$ws.Range("K73").Value = 10374.5
$ws.Range("M73").Value = -9438.5
$ws.Range("H97").Value = 1082.9166
$ws.Range("I97").Value = 1096.8334
$ws.Range("K97").Value = 1096.8334
$ws.Range("M97").Value = -600.8334
$ws.Range("H102").Value = 1130.625
$ws.Range("I102").Value = 930.2308
$ws.Range("J102").Value = 1999
$ws.Range("K102").Value = 930.2308
$ws.Range("L102").Value = 1999
$ws.Range("M102").Value = 691.7692
$ws.Range("N102").Value = -5243
$ws.Range("H107").Value = 804.4375
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()
$ws.Range("H113").Value = 2608.647
$ws.Range("I113").Value = 2621.6875
$ws.Range("K113").Value = 2621.6875
$ws.Range("M113").Value = -451.6875
$ws.Range("H122").Value = 5931.3335
$ws.Range("I122").Value = 6917.6
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 20752.8
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -18302.8
$ws.Range("N122").Value = -7900
$ws.Range("H124").Value = 72500
$ws.Range("J124").Value = 72500
$ws.Range("L124").Value = 72500
$ws.Range("N124").Value = -82320
$ws.Range("H133").Value = 64999.5
$ws.Range("J133").Value = 64999.5
$ws.Range("L133").Value = 64999.5
$ws.Range("N133").Value = -75119.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2980.8948
$ws.Range("I22").Value = 1158.8
$ws.Range("J22").Value = 3631.6428
$ws.Range("K22").Value = 1158.8
$ws.Range("L22").Value = 3631.6428
$ws.Range("M22").Value = -863.8
$ws.Range("N22").Value = -4221.6428
$ws.Range("H27").Value = 2980.8948
$ws.Range("I27").Value = 1158.8
$ws.Range("J27").Value = 3631.6428
$ws.Range("K27").Value = 1158.8
$ws.Range("L27").Value = 3631.6428
$ws.Range("M27").Value = -1051.8
$ws.Range("N27").Value = -3845.6428
$ws.Range("H40").Value = 3082.3076
$ws.Range("I40").Value = 2985.3684
$ws.Range("J40").Value = 3345.4285
$ws.Range("K40").Value = 2985.3684
$ws.Range("L40").Value = 3345.4285
$ws.Range("M40").Value = -2849.3684
$ws.Range("N40").Value = -3617.4285
$ws.Range("H68").Value = 4323.4
$ws.Range("I68").Value = 4841.75
$ws.Range("J68").Value = 2250
$ws.Range("K68").Value = 4841.75
$ws.Range("L68").Value = 2250
$ws.Range("M68").Value = -4092.75
$ws.Range("N68").Value = -3748
$ws.Range("H71").Value = 4323.4
$ws.Range("I71").Value = 4841.75
$ws.Range("J71").Value = 2250
$ws.Range("K71").Value = 24208.75
$ws.Range("L71").Value = 11250
$ws.Range("M71").Value = -20464.75
$ws.Range("N71").Value = -18738
$ws.Range("H74").Value = 25217
$ws.Range("J74").Value = 25217
$ws.Range("L74").Value = 25217
$ws.Range("N74").Value = -27213
$ws.Range("H77").Value = 25217
$ws.Range("J77").Value = 25217
$ws.Range("L77").Value = 75651
$ws.Range("N77").Value = -85635
$ws.Range("H122").Value = 6290.8
$ws.Range("I122").Value = 5851.3335
$ws.Range("J122").Value = 6950
$ws.Range("K122").Value = 17554.0005
$ws.Range("L122").Value = 20850
$ws.Range("M122").Value = -15104.0005
$ws.Range("N122").Value = -25750
$ws.Range("H124").Value = 99999
$ws.Range("J124").Value = 99999
$ws.Range("L124").Value = 99999
$ws.Range("N124").Value = -109819  # new cell
$ws.Range("H133").Value = 58326
$ws.Range("J133").Value = 58326
$ws.Range("L133").Value = 58326
$ws.Range("N133").Value = -63386
$ws.Range("H136").Value = 58826264
$ws.Range("I136").Value = 2450.1428
$ws.Range("K136").Value = 7350.428400000001
$ws.Range("M136").Value = -4800.428400000001
$ws.Range("H141").Value = 79998.5
$ws.Range("I141").Value = 79999
$ws.Range("J141").Value = 79998
$ws.Range("K141").Value = 79999
$ws.Range("L141").Value = 79998
$ws.Range("N141").Value = -90358
$ws.Range("M141").Value = -74819  # new cell

